# Update the "2024" sheet: a new "credit icici" SMS notification arrived,
# pushing the existing R47:S181 (message/time) rows down by one, and the
# "Broadband" category label moves from A189 down to A190.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Target values (after the shift) for columns R (message) and S (timestamp)
# for rows 47 through 181.
$data = @(
  @(47, 'credit icici', '2024-09-22 15:23:32'),
  @(48, 'credit icici', '2024-09-22 14:31:27'),
  @(49, 'balance your axis', '2024-09-22 08:02:25'),
  @(50, 'balance your axis', '2024-09-21 11:25:06'),
  @(51, 'balance your axis', '2024-09-21 10:34:04'),
  @(52, 'bal axis', '2024-09-21 07:56:12'),
  @(53, 'bal axisbank axis', '2024-09-21 07:50:18'),
  @(54, 'check the loan yo', '2024-09-20 15:37:11'),
  @(55, 'balance your axis', '2024-09-20 08:05:28'),
  @(56, 'bal axis', '2024-09-20 07:03:45'),
  @(57, 'axis', '2024-09-20 06:57:43'),
  @(58, 'dispute', '2024-09-19 22:46:00'),
  @(59, 'tamilnadu disclose it anyone', '2024-09-19 22:41:11'),
  @(60, 'dispute', '2024-09-19 22:33:39'),
  @(61, 'dispute', '2024-09-19 22:27:16'),
  @(62, 'your relationship', '2024-09-19 15:37:45'),
  @(63, 'value discovery debit icici', '2024-09-19 14:34:40'),
  @(64, 'debit', '2024-09-19 14:35:16'),
  @(65, 'balance your axis', '2024-09-19 11:05:17'),
  @(66, 'balance your axis', '2024-09-18 12:48:31'),
  @(67, 'your relationship', '2024-09-18 10:29:06'),
  @(68, 'balance your axis', '2024-09-18 10:28:28'),
  @(69, 'axis', '2024-09-18 08:12:44'),
  @(70, 'broker', '2024-09-18 04:09:58'),
  @(71, 'balance your axis', '2024-09-17 13:07:16'),
  @(72, 'dispute', '2024-09-16 12:53:44'),
  @(73, 'money google icici', '2024-09-16 12:53:29'),
  @(74, 'indusind', '2024-09-16 11:13:15'),
  @(75, 'communication feedback', '2024-09-16 11:13:15'),
  @(76, 'balance your axis', '2024-09-16 08:57:11'),
  @(77, 'balance your axis', '2024-09-16 07:57:00'),
  @(78, 'money google icici', '2024-09-15 21:06:00'),
  @(79, 'adani icici', '2024-09-15 13:10:50'),
  @(80, 'balance your axis', '2024-09-15 07:56:24'),
  @(81, 'bal axisbank w axis', '2024-09-15 07:12:01'),
  @(82, 'hdfc', '2024-09-14 21:25:23'),
  @(83, 'change the', '2024-09-12 21:16:38'),
  @(84, 'dispute', '2024-09-12 19:02:14'),
  @(85, 'congrats limit icici', '2024-09-12 19:03:39'),
  @(86, 'latest transaction pan', '2024-09-12 12:22:12'),
  @(87, 'assistance', '2024-09-12 12:17:33'),
  @(88, 'balance your axis', '2024-09-12 09:37:28'),
  @(89, 'bal axisbank', '2024-09-12 00:54:39'),
  @(90, 'your relationship', '2024-09-11 16:05:27'),
  @(91, 'login internet personal share', '2024-09-11 14:16:45'),
  @(92, 'balance your axis', '2024-09-11 12:45:33'),
  @(93, 'balance your axis', '2024-09-11 09:45:01'),
  @(94, 'axis', '2024-09-11 06:57:42'),
  @(95, 'money google icici', '2024-09-10 20:42:12'),
  @(96, 'dispute', '2024-09-10 20:42:34'),
  @(97, 'reward points cash', '2024-09-10 19:43:35'),
  @(98, 'balance your axis', '2024-09-10 13:32:42'),
  @(99, 'ach indianesign bal axisbank', '2024-09-10 13:22:37'),
  @(100, 'ach indianesign bal axisbank', '2024-09-10 13:22:37'),
  @(101, 'balance your axis', '2024-09-10 11:21:40'),
  @(102, 'your relationship', '2024-09-10 11:02:23'),
  @(103, 'bank bal broker', '2024-09-09 19:59:02'),
  @(104, 'beneficiary', '2024-09-09 15:48:10'),
  @(105, 'beneficiary saravanan', '2024-09-09 14:52:20'),
  @(106, 'bal axisbank', '2024-09-09 12:19:34'),
  @(107, 'bal axisbank', '2024-09-09 12:19:33'),
  @(108, 'dispute', '2024-09-09 12:17:30'),
  @(109, 'bal axisbank', '2024-09-09 12:04:31'),
  @(110, 'transfer freedom share anyone axis', '2024-09-09 11:56:19'),
  @(111, 'corporate internet share', '2024-09-09 11:40:49'),
  @(112, 'corporate internet share', '2024-09-09 11:39:30'),
  @(113, 'bal axisbank', '2024-09-09 11:38:16'),
  @(114, 'bal axisbank', '2024-09-09 11:38:16'),
  @(115, 'bal axisbank', '2024-09-09 11:38:15'),
  @(116, 'bal axisbank', '2024-09-09 11:38:15'),
  @(117, 'corporate internet share', '2024-09-09 11:35:34'),
  @(118, 'corporate internet share', '2024-09-09 11:32:23'),
  @(119, 'ift anbu tpar', '2024-09-09 11:27:52'),
  @(120, 'balance your axis', '2024-09-09 11:24:00'),
  @(121, 'corporate internet share', '2024-09-09 11:21:43'),
  @(122, 'corporate internet share', '2024-09-09 11:17:34'),
  @(123, 'corporate internet share', '2024-09-09 11:15:51'),
  @(124, 'corporate internet share', '2024-09-09 11:14:13'),
  @(125, 'anbu tparty bal axisbank', '2024-09-09 11:13:37'),
  @(126, 'corporate internet share', '2024-09-09 11:10:39'),
  @(127, 'corporate internet share', '2024-09-09 11:07:31'),
  @(128, 'corporate internet share', '2024-09-09 11:03:09'),
  @(129, 'saravanan', '2024-09-09 10:43:11'),
  @(130, 'balance your axis', '2024-09-09 08:10:16'),
  @(131, 'ekalaivan', '2024-09-08 18:40:34'),
  @(132, 'balance your axis', '2024-09-08 09:53:37'),
  @(133, 'balance your axis', '2024-09-07 12:12:22'),
  @(134, 'balance your axis', '2024-09-07 09:34:58'),
  @(135, 'bal axis', '2024-09-07 08:46:40'),
  @(136, 'axis', '2024-09-07 08:31:28'),
  @(137, 'your relationship', '2024-09-06 12:23:25'),
  @(138, 'balance your axis', '2024-09-06 09:55:31'),
  @(139, 'beneficiary', '2024-09-05 17:13:56'),
  @(140, 'coimbatore ramalinga', '2024-09-05 17:06:01'),
  @(141, 'beneficiary', '2024-09-05 17:04:10'),
  @(142, 'bal axisbank', '2024-09-05 16:52:25'),
  @(143, 'share anyone axis', '2024-09-05 16:38:59'),
  @(144, 'transfer anyone axis', '2024-09-05 16:35:58'),
  @(145, 'share anyone axis', '2024-09-05 16:31:34'),
  @(146, 'transfer', '2024-09-05 16:28:38'),
  @(147, 'bal axisbank axis', '2024-09-05 16:26:56'),
  @(148, 'bal axisbank', '2024-09-05 16:26:55'),
  @(149, 'transfer', '2024-09-05 16:25:07'),
  @(150, 'transfer', '2024-09-05 16:22:23'),
  @(151, 'share anyone axis', '2024-09-05 16:06:05'),
  @(152, 'internet bal axisbank', '2024-09-05 16:05:55'),
  @(153, 'transfer share anyone axis', '2024-09-05 16:03:14'),
  @(154, 'axis', '2024-09-05 15:57:15'),
  @(155, 'your net internet', '2024-09-05 15:57:15'),
  @(156, 'hear your feedback atm', '2024-09-05 14:21:08'),
  @(157, 'axis bna', '2024-09-05 14:18:32'),
  @(158, 'axis bna', '2024-09-05 14:13:16'),
  @(159, 'axis bna', '2024-09-05 14:15:23'),
  @(160, 'balance your axis', '2024-09-05 09:20:57'),
  @(161, 'bal axis', '2024-09-05 09:06:25'),
  @(162, 'broker', '2024-09-04 21:20:47'),
  @(163, 'exclusive on axis', '2024-09-04 13:21:05'),
  @(164, 'your corporate axis', '2024-09-04 11:46:10'),
  @(165, 'balance your axis', '2024-09-04 08:14:16'),
  @(166, 'axis', '2024-09-04 07:02:13'),
  @(167, 'bal axisbank w axis', '2024-09-04 06:53:15'),
  @(168, 'logging iob internet', '2024-09-03 20:09:12'),
  @(169, 'password internet', '2024-09-03 20:05:31'),
  @(170, 'logging iob internet', '2024-09-03 20:05:09'),
  @(171, 'internet', '2024-09-03 19:58:18'),
  @(172, 'login internet invalid', '2024-09-03 19:54:49'),
  @(173, 'login internet invalid', '2024-09-03 19:56:17'),
  @(174, 'corporate internet share', '2024-09-03 19:22:58'),
  @(175, 'login sbi internet personal do not share anyone', '2024-09-03 19:17:10'),
  @(176, 'login internet personal share', '2024-09-03 19:13:40'),
  @(177, 'internet verify it', '2024-09-03 19:05:49'),
  @(178, 'balance your axis', '2024-09-03 13:14:06'),
  @(179, 'lounge', '2024-09-03 13:08:08'),
  @(180, 'balance your axis', '2024-09-03 11:21:30'),
  @(181, 'broker', '2024-09-01 22:35:38')
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 18).Value = $row[1]
  $ws.Cells.Item($r, 19).Value = $row[2]
}

# "Broadband" label moves from row 189 to the newly appended row 190.
$ws.Cells.Item(189, 1).Value = ""
$ws.Cells.Item(190, 1).Value = "Broadband"
